$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("H54").Value = 16500.445
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 19786.285
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 19786.285
$ws.Range("M54").Value = -4514
$ws.Range("N54").Value = -20758.285
$ws.Range("H98").Value = 2976.25
$ws.Range("I98").Value = 3952.5
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 3952.5
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -2454.5
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 2976.25
$ws.Range("I122").Value = 3952.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 11857.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -9407.5
$ws.Range("N122").Value = -10900
$ws.Range("H125").Value = 7266.3335
$ws.Range("I125").Value = 700
$ws.Range("J125").Value = 7735.357
$ws.Range("K125").Value = 6300
$ws.Range("L125").Value = 69618.213
$ws.Range("M125").Value = -3840
$ws.Range("N125").Value = -74538.213
$ws.Range("H126").Value = 66733.336
$ws.Range("J126").Value = 66733.336
$ws.Range("L126").Value = 66733.336
$ws.Range("N126").Value = -76613.336
$ws.Range("H128").Value = 34890
$ws.Range("J128").Value = 34890
$ws.Range("L128").Value = 34890
$ws.Range("N128").Value = -44850
$ws.Range("H138").Value = 4917.5386
$ws.Range("I138").Value = 1253.4482
$ws.Range("J138").Value = 15543.4
$ws.Range("K138").Value = 3760.3446
$ws.Range("L138").Value = 46630.2
$ws.Range("M138").Value = 1379.6554
$ws.Range("N138").Value = -56910.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1840.3334
$ws.Range("I74").Value = 1816.6666
$ws.Range("K74").Value = 1816.6666
$ws.Range("M74").Value = -942.6666
$ws.Range("H77").Value = 1840.3334
$ws.Range("I77").Value = 1816.6666
$ws.Range("K77").Value = 9083.333000000001
$ws.Range("M77").Value = -4715.333000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 36360.5
$ws.Range("I76").Value = 25000
$ws.Range("J76").Value = 37983.43
$ws.Range("K76").Value = 25000
$ws.Range("L76").Value = 37983.43
$ws.Range("M76").Value = -24685
$ws.Range("N76").Value = -38613.43
$ws.Range("H79").Value = 36360.5
$ws.Range("I79").Value = 25000
$ws.Range("J79").Value = 37983.43
$ws.Range("K79").Value = 25000
$ws.Range("L79").Value = 37983.43
$ws.Range("M79").Value = -23908
$ws.Range("N79").Value = -40167.43
$ws.Range("H86").Value = 1635.4314
$ws.Range("I86").Value = 1638.4468
$ws.Range("K86").Value = 1638.4468
$ws.Range("M86").Value = -515.4467999999999
$ws.Range("H89").Value = 1635.4314
$ws.Range("I89").Value = 1638.4468
$ws.Range("K89").Value = 8192.234
$ws.Range("M89").Value = -2576.234

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2714.1562
$ws.Range("I31").Value = 1417
$ws.Range("K31").Value = 1417
$ws.Range("M31").Value = -1122
$ws.Range("H34").Value = 2714.1562
$ws.Range("I34").Value = 1417
$ws.Range("K34").Value = 1417
$ws.Range("M34").Value = -1215
$ws.Range("H122").Value = 7707.84
$ws.Range("I122").Value = 5021.467
$ws.Range("K122").Value = 15064.401
$ws.Range("M122").Value = -12614.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 12828874
$ws.Range("I5").Value = 753.75
$ws.Range("J5").Value = 18530260
$ws.Range("K5").Value = 2261.25
$ws.Range("L5").Value = 55590780
$ws.Range("M5").Value = -2149.25
$ws.Range("N5").Value = -55591004
$ws.Range("H14").Value = 1451.6666
$ws.Range("I14").Value = 1451.6666
$ws.Range("K14").Value = 4354.9998
$ws.Range("M14").Value = -4181.9998
$ws.Range("H122").Value = 662.24243
$ws.Range("I122").Value = 305.3125
$ws.Range("K122").Value = 2747.8125
$ws.Range("M122").Value = -297.8125
$ws.Range("H127").Value = 1795.3793
$ws.Range("J127").Value = 1795.3793
$ws.Range("L127").Value = 5386.1379
$ws.Range("N127").Value = -15306.1379
$ws.Range("H131").Value = 1360.9667
$ws.Range("J131").Value = 1154.9375
$ws.Range("L131").Value = 3464.8125
$ws.Range("N131").Value = -13544.8125
$ws.Range("H135").Value = 12828874
$ws.Range("I135").Value = 753.75
$ws.Range("J135").Value = 18530260
$ws.Range("K135").Value = 6783.75
$ws.Range("L135").Value = 166772340
$ws.Range("M135").Value = -4248.75
$ws.Range("N135").Value = -166777410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 12200
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 15000
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -673
$ws.Range("N55").Value = -15654
$ws.Range("H126").Value = 2799.9333
$ws.Range("J126").Value = 3090.818
$ws.Range("L126").Value = 9272.454000000002
$ws.Range("N126").Value = -14212.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3363
$ws.Range("I68").Value = 2602
$ws.Range("J68").Value = 3616.6667
$ws.Range("K68").Value = 2602
$ws.Range("L68").Value = 3616.6667
$ws.Range("M68").Value = -1853
$ws.Range("N68").Value = -5114.6667
$ws.Range("H71").Value = 3363
$ws.Range("I71").Value = 2602
$ws.Range("J71").Value = 3616.6667
$ws.Range("K71").Value = 13010
$ws.Range("L71").Value = 18083.3335
$ws.Range("M71").Value = -9266

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3888
$ws.Range("J62").Value = 3850.6667
$ws.Range("L62").Value = 3850.6667
$ws.Range("N62").Value = -5098.6667
$ws.Range("H65").Value = 3888
$ws.Range("J65").Value = 3850.6667
$ws.Range("L65").Value = 19253.3335
$ws.Range("N65").Value = -25493.3335
$ws.Range("H107").Value = 3073.0715
$ws.Range("I107").Value = 601.7
$ws.Range("K107").Value = 1805.1
$ws.Range("M107").Value = 114.8999999999999
